# Populate Sheet1 column A with an NFL "Team Name" header + all 32 team
# names (ESPN NFL test data), then reproduce the author's later edit of
# inserting the two AFC South teams that were originally missing
# (Indianapolis Colts, Tennessee Titans) as new rows, which is why they
# land at the end of the shared-string table while appearing mid-list in
# the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$teams = @(
    "Team Name",
    "Buffalo Bills",
    "Miami Dolphins",
    "New England Patriots",
    "New York Jets",
    "Baltimore Ravens",
    "Cincinnati Bengals",
    "Cleveland Browns",
    "Pittsburgh Steelers",
    "Houston Texans",
    "Jacksonville Jaguars",
    "Denver Broncos",
    "Kansas City Chiefs",
    "Las Vegas Raiders",
    "Los Angeles Chargers",
    "Dallas Cowboys",
    "New York Giants",
    "Philadelphia Eagles",
    "Washington Commanders",
    "Chicago Bears",
    "Detroit Lions",
    "Green Bay Packers",
    "Minnesota Vikings",
    "Atlanta Falcons",
    "Carolina Panthers",
    "New Orleans Saints",
    "Tampa Bay Buccaneers",
    "Arizona Cardinals",
    "Los Angeles Rams",
    "San Francisco 49ers",
    "Seattle Seahawks"
)

for ($i = 0; $i -lt $teams.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $teams[$i]
}

# Insert "Indianapolis Colts" as a new row 11 (AFC South), shifting rows 11+ down
$ws.Rows(11).Insert()
$ws.Cells.Item(11, 1).Value = "Indianapolis Colts"

# Insert "Tennessee Titans" as a new row 13 (AFC South), shifting rows 13+ down
$ws.Rows(13).Insert()
$ws.Cells.Item(13, 1).Value = "Tennessee Titans"

# Restore the view: zoomed in, scrolled down, cursor left on L23
$excel.ActiveWindow.Zoom = 175
$ws.Range("L23").Select()

